$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove two outdated books at the top of the list ---
# Row 6: "Povesti in romana si germana" (Sojka Anna)
# Row 7: "Kathie si hipopotamul" (Vargas Llosa Mario)
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()

# After the deletions above, the remaining books shift up:
#   6 Venus ia foc Marte e de gheata
#   7 Cu ochii mintii
#   8 In apararea pietelor
#   9 Adobe InDesign CC. ...
#  10 Щенок Элфи или Не хочу быть один! <- remove this one too
$ws.Rows.Item(10).Delete()

# --- Insert a new book row before "Adobe InDesign CC. ..." (now row 9) ---
$ws.Rows.Item(9).Insert()

# Clone the formatting (borders / wrap / alignment) of the Adobe row (now
# pushed down to row 10) onto the freshly inserted, blank row 9.
$ws.Range("B10:I10").Copy()
$ws.Range("B9:I9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# This is a long title, so it needs the taller, word-wrapped row like the
# other long titles (Adobe / the removed Russian book) already use.
$ws.Rows.Item(9).RowHeight = 65
$ws.Cells.Item(9, 2).WrapText = $true

# --- Fill in the new book's data ---
$ws.Cells.Item(9, 2).Value = "Cele 12 elemente ale managementului performant"
$ws.Cells.Item(9, 3).Value = "Wagner R."
$ws.Cells.Item(9, 4).Value = "All"
$ws.Cells.Item(9, 5).Value = "'2009"
$ws.Cells.Item(9, 6).Value = "'272"
$ws.Cells.Item(9, 7).Value = "Management"
$ws.Cells.Item(9, 8).Value = "'9789737241887"
$ws.Cells.Item(9, 9).Value = "'155"

# --- Refresh the export timestamp footer to the current local time ---
$ws.Cells.Item(13, 2).Value = "Data exportării: 05/22/2022 14:05"
